$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.737.70'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.863.46'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.038'
$ws.Range("E4").Value = '  +0.94%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.43'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4423'
$ws.Range("E7").Value = '  +1.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3809'
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07465'
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8870'
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.72'
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D12").Value = '1.875.78'
$ws.Range("E12").Value = '  -5.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.550'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.778'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07210'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.31'
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.039'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009144'
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.56'
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").Value = '27.749.38'
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.316'
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.30'
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").Value = '2.092.04'
$ws.Range("E24").Value = '  -3.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.020'
$ws.Range("E25").Value = '  +6.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.28'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.366'
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.986'
$ws.Range("E29").Value = '  +3.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.07'
$ws.Range("E30").Value = '  +3.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09057'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.229'
$ws.Range("E32").Value = '  +2.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7805'
$ws.Range("E33").Value = '  +2.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.028'
$ws.Range("E34").Value = '  +5.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.599'
$ws.Range("E35").Value = '  +2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.034'
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01988'
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05362'
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.885'
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5215'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1696'
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.908'
$ws.Range("E43").Value = '  +5.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.730'
$ws.Range("E44").Value = '  +3.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.99'
$ws.Range("E45").Value = '  +2.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.71'
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06686'
$ws.Range("E47").Value = '  +6.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.036'
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.714'
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4734'
$ws.Range("E50").Value = '  +2.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.916'
$ws.Range("E51").Value = '  +1.95%  '
